# Fix Training Data Issue
# The BF column ("Date") held malformed strings like "4-28-2007-08" that were
# actually meant to read the proper ISO-ish date "2008-04-28" (NBA stats for
# 2007-08 season, as of 4/28/2008). Correct all 30 data rows (BF2:BF31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 58)  # column BF = 58
    if ($cell.Text -eq "4-28-2007-08") {
        # Force the assignment to stay plain text instead of being
        # auto-parsed into a date serial number, then drop the temporary
        # number-format override so the cell keeps its original (default)
        # style.
        $cell.NumberFormat = "@"
        $cell.Value = "2008-04-28"
        $cell.ClearFormats()
    }
}
